$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change cell B11 value from "R40" to the string "1"
$ws.Range("B11").Value = "1"
